$wb = $excel.ActiveWorkbook

# "Generate Report for Handoff" - update the per-language handoff tracking
# sheets with the new handoff generation: Priority = "ht" and a refreshed
# "Latest Handoff Datetime" for every source file row.

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E2").Value = "ht"
$zhcn.Range("H2").Value = "2016-11-29 07:14:22"
$zhcn.Range("E3").Value = "ht"
$zhcn.Range("H3").Value = "2016-11-29 07:14:22"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E2").Value = "ht"
$dede.Range("H2").Value = "2016-11-29 07:14:35"
$dede.Range("E3").Value = "ht"
$dede.Range("H3").Value = "2016-11-29 07:14:35"

# The de-de handoff datetime is also mirrored on the Overview summary sheet
# (it shares the same underlying "Latest HO Xliff Generate Date" value), so
# refresh it there too.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G2").Value = "2016-11-29 07:14:35"
$overview.Range("G3").Value = "2016-11-29 07:14:35"
